$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Error handling" task row is being moved from row 6 up to row 4
# (inserted right after "Make the path from Max->Model Viewer seamless"
# and before "Replace XNAMATH with a Rorn maths library"), and its text is
# being extended in the process. The rows that used to sit at 4 and 5 each
# shift down by one to make room; everything below row 6 is unaffected.
#
# Before (rows 2-9):
#   2 Make exporter a GUP.  Build UI & hook data into the max files    | 14
#   3 Make the path from Max->Model Viewer seamless                    | 4
#   4 Replace XNAMATH with a Rorn maths library                        | 14
#   5 Refactoring - we need consistency across the board                | 21
#   6 Error handling - go on a robustness run, set standards for
#     future work                                                       | 21
#   7 Textured surfaces                                                 | 35
#   8 Vertex welding in model compiler                                  | 14
#   9 Complete the Rorn Maths library                                   | 35
#
# After (rows 2-9):
#   2 Make exporter a GUP.  Build UI & hook data into the max files    | 14
#   3 Make the path from Max->Model Viewer seamless                    | 4
#   4 Error handling - go on a robustness run, set standards for
#     future work.  Ensure anythign that can go wrong is handled.       | 21
#   5 Replace XNAMATH with a Rorn maths library                        | 14
#   6 Refactoring - we need consistency across the board                | 21
#   7 Textured surfaces                                                 | 35
#   8 Vertex welding in model compiler                                  | 14
#   9 Complete the Rorn Maths library                                   | 35

# Snapshot the old values first so shifting doesn't clobber data we still
# need to read.
$task4 = $ws.Range("A4").Value()
$est4  = $ws.Range("B4").Value()
$task5 = $ws.Range("A5").Value()
$est5  = $ws.Range("B5").Value()

# Shift old row 4 -> row 5, old row 5 -> row 6.
$ws.Range("A5").Value = $task4
$ws.Range("B5").Value = $est4

$ws.Range("A6").Value = $task5
$ws.Range("B6").Value = $est5

# Put the updated "Error handling" task into row 4; estimate is unchanged (21).
$ws.Range("A4").Value = "Error handling - go on a robustness run, set standards for future work.  Ensure anythign that can go wrong is handled."
$ws.Range("B4").Value = 21

# Match the cosmetic view-state changes captured in the diff: selection
# moved to A5, and the first column's width grew (bestFit-style) to
# accommodate the new, longer "Error handling" text.
$ws.Range("A5").Select()
$ws.Columns("A:A").ColumnWidth = 105.66666666666667
